$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the styling of the
# existing header cells (bold font, border, centered alignment - style "1").
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Column I (I0) values for rows 2-66.
$iValues = @(9,8,9,8,8,8,9,7,7,8,5,8,7,8,8,6,9,9,8,8,8,8,8,7,7,7,7,9,8,7,8,8,9,5,7,6,6,6,8,8,8,6,8,9,8,7,8,5,9,7,8,8,7,9,8,7,9,8,7,6,4,6,4,6,3)

# Column J (IF) values for rows 2-66.
$jValues = @(9,8,9,8,8,8,9,8,7,8,6,8,8,8,8,7,9,9,8,8,8,8,8,7,7,7,7,9,8,7,8,8,9,6,7,6,7,7,9,8,8,6,8,9,8,7,8,6,9,8,8,8,7,9,8,7,9,8,7,6,5,6,4,6,3)

for ($n = 0; $n -lt $iValues.Length; $n++) {
    $row = $n + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$n]
    $ws.Cells.Item($row, 10).Value = $jValues[$n]
}
